$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price (D) and 1h volume-change (E) values for rows 2-51.
# Some prices are plain decimals (e.g. "1.00", "257.41") that Excel would
# otherwise auto-convert to numbers, stripping meaningful trailing zeros;
# those cells are forced to text format before the value is written so the
# literal string is preserved, matching how the source data is stored.
# Rows whose Price column did not change keep $null in D and are skipped.
$updates = @(
    @{ Row = 2; D = "43.491.69"; E = "  +2.43%  "; ForceText = $False },
    @{ Row = 3; D = "2.178.31"; E = "  -0.25%  "; ForceText = $False },
    @{ Row = 4; D = $null; E = "  +0.14%  "; ForceText = $False },
    @{ Row = 5; D = "257.41"; E = "  +1.08%  "; ForceText = $True },
    @{ Row = 6; D = "79.93"; E = "  +8.77%  "; ForceText = $True },
    @{ Row = 7; D = "0.620"; E = "  +2.17%  "; ForceText = $True },
    @{ Row = 8; D = $null; E = "  +0.07%  "; ForceText = $False },
    @{ Row = 9; D = "0.586"; E = "  +0.92%  "; ForceText = $True },
    @{ Row = 10; D = "42.60"; E = "  +5.15%  "; ForceText = $True },
    @{ Row = 11; D = "0.0912"; E = "  -0.53%  "; ForceText = $True },
    @{ Row = 12; D = $null; E = "  +2.55%  "; ForceText = $False },
    @{ Row = 13; D = "6.90"; E = "  +1.75%  "; ForceText = $True },
    @{ Row = 14; D = "2.508.99"; E = "  +0.00%  "; ForceText = $False },
    @{ Row = 15; D = "14.16"; E = "  -0.06%  "; ForceText = $True },
    @{ Row = 16; D = "2.181.97"; E = "  +0.27%  "; ForceText = $False },
    @{ Row = 17; D = "0.769"; E = "  -0.53%  "; ForceText = $True },
    @{ Row = 18; D = "43.445.10"; E = "  +2.55%  "; ForceText = $False },
    @{ Row = 19; D = $null; E = "  -0.27%  "; ForceText = $False },
    @{ Row = 20; D = "69.63"; E = "  -1.50%  "; ForceText = $True },
    @{ Row = 21; D = "5.86"; E = "  -0.19%  "; ForceText = $True },
    @{ Row = 22; D = "2.37"; E = "  +11.27%  "; ForceText = $True },
    @{ Row = 23; D = "229.03"; E = "  +0.89%  "; ForceText = $True },
    @{ Row = 24; D = "8.76"; E = "  -6.99%  "; ForceText = $True },
    @{ Row = 25; D = $null; E = "  +0.13%  "; ForceText = $False },
    @{ Row = 26; D = "42.16"; E = "  +14.01%  "; ForceText = $True },
    @{ Row = 27; D = "10.58"; E = "  +0.88%  "; ForceText = $True },
    @{ Row = 28; D = $null; E = "  -0.31%  "; ForceText = $False },
    @{ Row = 29; D = $null; E = "  +4.98%  "; ForceText = $False },
    @{ Row = 30; D = $null; E = "  +1.62%  "; ForceText = $False },
    @{ Row = 31; D = "172.90"; E = "  +1.22%  "; ForceText = $True },
    @{ Row = 32; D = "20.24"; E = "  +1.06%  "; ForceText = $True },
    @{ Row = 33; D = "0.0863"; E = "  +6.98%  "; ForceText = $True },
    @{ Row = 34; D = "5.25"; E = "  +2.39%  "; ForceText = $True },
    @{ Row = 35; D = $null; E = "  +4.14%  "; ForceText = $False },
    @{ Row = 36; D = $null; E = "  +0.86%  "; ForceText = $False },
    @{ Row = 37; D = "4.42"; E = "  +3.93%  "; ForceText = $True },
    @{ Row = 38; D = "0.0348"; E = "  +2.76%  "; ForceText = $True },
    @{ Row = 39; D = "12.92"; E = "  +9.19%  "; ForceText = $True },
    @{ Row = 40; D = "2.80"; E = "  +13.24%  "; ForceText = $True },
    @{ Row = 41; D = "2.08"; E = "  +1.01%  "; ForceText = $True },
    @{ Row = 42; D = "62.17"; E = "  +4.66%  "; ForceText = $True },
    @{ Row = 43; D = "5.40"; E = "  +4.91%  "; ForceText = $True },
    @{ Row = 44; D = "0.197"; E = "  -0.02%  "; ForceText = $True },
    @{ Row = 45; D = "100.08"; E = "  -2.24%  "; ForceText = $True },
    @{ Row = 46; D = "0.0977"; E = "  +0.61%  "; ForceText = $True },
    @{ Row = 47; D = "8.15"; E = "  -1.00%  "; ForceText = $True },
    @{ Row = 48; D = $null; E = "  +3.62%  "; ForceText = $False },
    @{ Row = 49; D = $null; E = "  +0.64%  "; ForceText = $False },
    @{ Row = 50; D = "1.52"; E = "  +25.01%  "; ForceText = $True },
    @{ Row = 51; D = "0.434"; E = "  -6.41%  "; ForceText = $True }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        if ($u.ForceText) {
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

Write-Host "Updated $($updates.Count) rows"
